$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell H2: fix wording ("wether" -> "whether", "Contct" -> "Contact")
#     and the "Taining" bold run ("aining" -> "ainting") ---
$cell = $ws.Range("H2")

$newText = "1. Open the CHHS URL.`n2. Check whether header tab About Us is clickable or not.`n3  Check whether header tab For Families is clickable or not.`n4. Check whether header tab Tainting is clickable or not.`n5. Check whether header tab Contact Us is clickable or not.`n6. Check whether header tab Search is clickable or not.`n7. Check whether header tab Login is clickable or not.`n8. Check whether header tab Help is clickable or not.`n9. Check whether Family Registration link is clickable or not."

$cell.Value2 = $newText

# Re-apply the bold formatting to the same words/phrases as before the edit
$cell.Characters(13, 4).Font.Bold = $true    # CHHS
$cell.Characters(51, 9).Font.Bold = $true    # About Us
$cell.Characters(109, 12).Font.Bold = $true  # For Families
$cell.Characters(172, 7).Font.Bold = $true   # ainting (part of Tainting)
$cell.Characters(229, 10).Font.Bold = $true  # Contact Us
$cell.Characters(289, 6).Font.Bold = $true   # Search
$cell.Characters(345, 5).Font.Bold = $true   # Login
$cell.Characters(400, 4).Font.Bold = $true   # Help
$cell.Characters(443, 19).Font.Bold = $true  # Family Registration

# --- Move the active selection to H3 (as recorded in the saved view state) ---
$ws.Range("H3").Select()
